$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 258; existing rows 258-333 shift down to 259-334.
$ws.Range("A258").EntireRow.Insert()

# Populate the newly inserted row 258 with the new record's data.
$ws.Range("A258").Value = 10
$ws.Range("B258").Value = 'Vega Modelo de Temuco'
$ws.Range("C258").Value = 'La Araucanía'
$ws.Range("D258").Value = 44876
$ws.Range("E258").Value = 9
$ws.Range("F258").Value = 'Fruta'
$ws.Range("G258").Value = 100102
$ws.Range("H258").Value = 'Cítricos'
$ws.Range("I258").Value = 100102006
$ws.Range("J258").Value = 'Pomelo'
$ws.Range("K258").Value = 'Start Ruby'
$ws.Range("L258").Value = 'Primera'
$ws.Range("M258").Value = 80
$ws.Range("N258").Value = 13000
$ws.Range("O258").Value = 13000
$ws.Range("P258").Value = 13000
$ws.Range("Q258").Value = '$/bandeja 15 kilos granel'
$ws.Range("R258").Value = "Región de O'Higgins"
$ws.Range("S258").Value = 867
$ws.Range("T258").Value = 15
